# Weekly update: insert the newest price observation as a new row right
# after the existing row 38 (row 39), pushing all the subsequent
# historical rows down by one. The new row 147 is therefore the old
# row 146's data, shifted down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 39 - shifts rows 39:146 down to 40:147
# and extends the used range to A1:R147 (Excel also copies row 38's
# formatting, e.g. the date number format on column D, onto the new row).
$ws.Rows.Item(39).EntireRow.Insert()

# Populate the newly inserted row 39 with this week's new observation.
$ws.Cells.Item(39, 1).Value = 10
$ws.Cells.Item(39, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(39, 3).Value = "La Araucanía"
$ws.Cells.Item(39, 4).Value = (Get-Date -Year 2022 -Month 5 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(39, 5).Value = 9
$ws.Cells.Item(39, 6).Value = 100112012
$ws.Cells.Item(39, 7).Value = "Espinaca"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 30
$ws.Cells.Item(39, 11).Value = 9000
$ws.Cells.Item(39, 12).Value = 9000
$ws.Cells.Item(39, 13).Value = 9000
$ws.Cells.Item(39, 14).Value = "`$/docena de atados"
$ws.Cells.Item(39, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(39, 16).Value = 3000
$ws.Cells.Item(39, 17).Value = 3
$ws.Cells.Item(39, 18).Value = "Hortaliza"
